$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 8500
$ws.Range("I34").Value = 2750
$ws.Range("J34").Value = 20000
$ws.Range("K34").Value = 2750
$ws.Range("L34").Value = 20000
$ws.Range("M34").Value = -2547
$ws.Range("N34").Value = -20406

$ws.Range("H36").Value = 8500
$ws.Range("I36").Value = 2750
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 2750
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = -2035
$ws.Range("N36").Value = -21430

$ws.Range("H41").Value = 120.125
$ws.Range("I41").Value = 150
$ws.Range("J41").Value = 102.2
$ws.Range("K41").Value = 150
$ws.Range("L41").Value = 102.2
$ws.Range("M41").Value = 290
$ws.Range("N41").Value = -982.2

$ws.Range("H42").Value = 147.07692
$ws.Range("I42").Value = 22.4
$ws.Range("J42").Value = 225
$ws.Range("K42").Value = 67.19999999999999
$ws.Range("L42").Value = 675
$ws.Range("M42").Value = 162.8
$ws.Range("N42").Value = -1135

$ws.Range("H53").Value = 214.1
$ws.Range("I53").Value = 220.14285
$ws.Range("J53").Value = 200
$ws.Range("K53").Value = 220.14285
$ws.Range("L53").Value = 200
$ws.Range("M53").Value = 416.85715
$ws.Range("N53").Value = -1474

$ws.Range("H107").Value = 4125.5
$ws.Range("I107").Value = 6252.5
$ws.Range("J107").Value = 1998.5
$ws.Range("K107").Value = 6252.5
$ws.Range("L107").Value = 1998.5
$ws.Range("M107").Value = -4332.5
$ws.Range("N107").Value = -5838.5

$ws.Range("H132").Value = 3030.2126
$ws.Range("I132").Value = 2878.7173
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 8636.151899999999
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -6106.151899999999
$ws.Range("N132").Value = -35057

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H135").Value = 2778.0715
$ws.Range("I135").Value = 2530.2307
$ws.Range("J135").Value = 6000
$ws.Range("K135").Value = 22772.0763
$ws.Range("L135").Value = 54000
$ws.Range("M135").Value = -20237.0763
$ws.Range("N135").Value = -59070

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws.Range("H137").Value = 3058.1714
$ws.Range("I137").Value = 2661.7407
$ws.Range("J137").Value = 4396.125
$ws.Range("K137").Value = 7985.222099999999
$ws.Range("L137").Value = 13188.375
$ws.Range("M137").Value = -5435.222099999999
$ws.Range("N137").Value = -18288.375

$ws.Range("H138").Value = 2272.5518
$ws.Range("I138").Value = 3127.7144
$ws.Range("J138").Value = 2000.4546
$ws.Range("K138").Value = 9383.143199999999
$ws.Range("L138").Value = 6001.3638
$ws.Range("M138").Value = -4243.143199999999
$ws.Range("N138").Value = -16281.3638

$ws.Range("H141").Value = 5150.1113
$ws.Range("I141").Value = 2005.5555
$ws.Range("J141").Value = 8294.666999999999
$ws.Range("K141").Value = 6016.666499999999
$ws.Range("L141").Value = 24884.001
$ws.Range("M141").Value = -836.6664999999994
$ws.Range("N141").Value = -35244.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2411.0881
$ws.Range("J61").Value = 3340.9333
$ws.Range("L61").Value = 3340.9333
$ws.Range("N61").Value = -3764.9333

$ws.Range("H63").Value = 6128.4
$ws.Range("I63").Value = 3499.3333
$ws.Range("J63").Value = 7255.143
$ws.Range("K63").Value = 3499.3333
$ws.Range("L63").Value = 7255.143
$ws.Range("M63").Value = -2813.3333
$ws.Range("N63").Value = -8627.143

$ws.Range("H66").Value = 6128.4
$ws.Range("I66").Value = 3499.3333
$ws.Range("J66").Value = 7255.143
$ws.Range("K66").Value = 17496.6665
$ws.Range("L66").Value = 36275.715
$ws.Range("M66").Value = -14064.6665
$ws.Range("N66").Value = -43139.715

$ws.Range("H132").Value = 3125.258
$ws.Range("I132").Value = 2860.9614
$ws.Range("J132").Value = 4499.6
$ws.Range("K132").Value = 8582.8842
$ws.Range("L132").Value = 13498.8
$ws.Range("M132").Value = -6052.8842
$ws.Range("N132").Value = -18558.8

$ws.Range("H136").Value = 2411.0881
$ws.Range("J136").Value = 3340.9333
$ws.Range("L136").Value = 10022.7999
$ws.Range("N136").Value = -15122.7999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 24999
$ws.Range("J35").Value = 24999
$ws.Range("L35").Value = 24999
$ws.Range("N35").Value = -25619

$ws.Range("H82").Value = 16201.15
$ws.Range("I82").Value = 3723.1
$ws.Range("K82").Value = 3723.1
$ws.Range("M82").Value = -3340.1

$ws.Range("H85").Value = 16201.15
$ws.Range("I85").Value = 3723.1
$ws.Range("K85").Value = 3723.1
$ws.Range("M85").Value = -2397.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1494.8334
$ws.Range("I5").Value = 446.66666
$ws.Range("J5").Value = 1844.2222
$ws.Range("K5").Value = 446.66666
$ws.Range("L5").Value = 1844.2222
$ws.Range("M5").Value = -334.66666
$ws.Range("N5").Value = -2068.2222

$ws.Range("H29").Value = 4000
$ws.Range("J29").Value = 4000
$ws.Range("L29").Value = 4000
$ws.Range("N29").Value = -4586

$ws.Range("H132").Value = 16669037
$ws.Range("I132").Value = 1671.7142
$ws.Range("J132").Value = 55559556
$ws.Range("K132").Value = 5015.142599999999
$ws.Range("L132").Value = 166678668
$ws.Range("M132").Value = -2485.142599999999
$ws.Range("N132").Value = -166683728

$ws.Range("H134").Value = 1630
$ws.Range("J134").Value = 2000
$ws.Range("L134").Value = 6000
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 458
$ws.Range("I2").Value = 63.333332
$ws.Range("J2").Value = 1050
$ws.Range("K2").Value = 63.333332
$ws.Range("L2").Value = 1050
$ws.Range("M2").Value = 49.666668
$ws.Range("N2").Value = -1276

$ws.Range("H108").Value = 32000
$ws.Range("J108").Value = 32000
$ws.Range("L108").Value = 32000
$ws.Range("N108").Value = -39680

$ws.Range("H127").Value = 79326
$ws.Range("J127").Value = 79326
$ws.Range("L127").Value = 79326
$ws.Range("N127").Value = -89246

$ws.Range("H132").Value = 2529.5
$ws.Range("I132").Value = 1738.1538
$ws.Range("J132").Value = 3464.7273
$ws.Range("K132").Value = 5214.4614
$ws.Range("L132").Value = 10394.1819
$ws.Range("M132").Value = -2684.4614
$ws.Range("N132").Value = -15454.1819

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10300
$ws.Range("J22").Value = 13244.333
$ws.Range("L22").Value = 13244.333
$ws.Range("N22").Value = -13834.333

$ws.Range("H27").Value = 10300
$ws.Range("J27").Value = 13244.333
$ws.Range("L27").Value = 13244.333
$ws.Range("N27").Value = -13458.333

$ws.Range("H46").Value = 1099
$ws.Range("I46").Value = 1215
$ws.Range("J46").Value = 925
$ws.Range("K46").Value = 1215
$ws.Range("L46").Value = 925
$ws.Range("M46").Value = -1027
$ws.Range("N46").Value = -1301

$ws.Range("H55").Value = 388.71875
$ws.Range("I55").Value = 223.38461
$ws.Range("J55").Value = 501.8421
$ws.Range("K55").Value = 223.38461
$ws.Range("L55").Value = 501.8421
$ws.Range("M55").Value = -50.38461000000001
$ws.Range("N55").Value = -847.8421000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1377.4517
$ws.Range("I126").Value = 1091.7222
$ws.Range("J126").Value = 1773.0769
$ws.Range("K126").Value = 3275.1666
$ws.Range("L126").Value = 5319.2307
$ws.Range("M126").Value = -805.1665999999996
$ws.Range("N126").Value = -10259.2307

$ws.Range("H136").Value = 2246.5715
$ws.Range("I136").Value = 1772.8334
$ws.Range("J136").Value = 5089
$ws.Range("K136").Value = 5318.5002
$ws.Range("L136").Value = 15267
$ws.Range("M136").Value = -2768.5002
$ws.Range("N136").Value = -20367
